# Weekly update: prepend a new "Cilantro" price report (Vega Monumental
# Concepción) for the latest date, pushing all existing rows down by two
# rows (the sheet stores data as Primera/Segunda row-pairs per date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (row 78 is the first
# data row right after the header row). Everything currently at rows
# 78:215 shifts down to 80:217, and the sheet's dimension grows from
# R215 to R217 automatically.
$ws.Rows("78:79").Insert()

# New row 78: "Primera" quality entry for the newest date.
$ws.Range("A78").Value = 11
$ws.Range("B78").Value = "Vega Monumental Concepción"
$ws.Range("C78").Value = "Bíobío"
$ws.Range("D78").Value = 44771
$ws.Range("E78").Value = 8
$ws.Range("F78").Value = 100112040
$ws.Range("G78").Value = "Cilantro"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 200
$ws.Range("K78").Value = 700
$ws.Range("L78").Value = 800
$ws.Range("M78").Value = 750
$ws.Range("N78").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O78").Value = "Región de Ñuble"
$ws.Range("P78").Value = 750
$ws.Range("Q78").Value = 1
$ws.Range("R78").Value = "Hortaliza"

# New row 79: "Segunda" quality entry for the same date.
$ws.Range("A79").Value = 11
$ws.Range("B79").Value = "Vega Monumental Concepción"
$ws.Range("C79").Value = "Bíobío"
$ws.Range("D79").Value = 44771
$ws.Range("E79").Value = 8
$ws.Range("F79").Value = 100112040
$ws.Range("G79").Value = "Cilantro"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Segunda"
$ws.Range("J79").Value = 100
$ws.Range("K79").Value = 600
$ws.Range("L79").Value = 600
$ws.Range("M79").Value = 600
$ws.Range("N79").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O79").Value = "Región de Ñuble"
$ws.Range("P79").Value = 600
$ws.Range("Q79").Value = 1
$ws.Range("R79").Value = "Hortaliza"
